$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-14 (old dimension A1:N14 -> new dimension A1:N9)
$ws.Range("A10:N14").Delete()

# New data for rows 4-9 (re-run of the scheduling pre-test covering more zones)
$data = @(
    @(251346, "R3", 36.5, 70.16363636363636, "2025-04-10 10:41:00", "2025-04-10 11:17:30", "2025-04-10 11:17:30", "2025-04-10 12:27:39", 3859,  "bobina", "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9", 8,  70, 4),
    @(250866, "R3", 102,  104.4081632653061, "2025-04-10 07:18:00", "2025-04-10 09:00:00", "2025-04-10 09:00:00", "2025-04-10 10:44:24", 5116,  "bobina", "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9", 6,  70, 4),
    @(251550, "R3", 50,   727.5714285714286, "2025-04-10 10:44:24", "2025-04-10 11:34:24", "2025-04-10 11:34:24", "2025-04-14 07:41:58", 35651, "bobina", "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9", 3,  70, 4),
    @(251109, "R6", 112,  266.5915492957747, "2025-04-10 13:25:00", "2025-04-11 07:17:00", "2025-04-11 07:17:00", "2025-04-11 11:43:35", 18928, "bobina", "R6",                                                16, 70, 4),
    @(251088, "R9", 35,   89.6376811594203,  "2025-04-10 07:22:00", "2025-04-10 07:57:00", "2025-04-10 07:57:00", "2025-04-10 09:26:38", 6185,  "bobina", "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9", 3,  70, 4),
    @(235572, "R9", 35,   144.3188405797102, "2025-04-10 09:26:38", "2025-04-10 10:01:38", "2025-04-10 10:01:38", "2025-04-10 12:25:57", 9958,  "bobina", "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9",          5,  70, 4)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $r++
}
